$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.693.76'
$ws.Range("E2").Value = '  +1.90%  '
$ws.Range("D3").Value = '3.164.85'
$ws.Range("E3").Value = '  +1.53%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '529.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.50%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.537'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +14.31%  '
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.438'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.111'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.96%  '
$ws.Range("E12").Value = '  +2.56%  '
$ws.Range("D13").Value = '3.709.32'
$ws.Range("E13").Value = '  +1.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("E15").Value = '  +3.78%  '
$ws.Range("D16").Value = '58.714.13'
$ws.Range("E16").Value = '  +1.71%  '
$ws.Range("D17").Value = '3.174.41'
$ws.Range("E17").Value = '  +1.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.98'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '376.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.27%  '
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("E23").Value = '  +4.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.16%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.67%  '
$ws.Range("D28").Value = '0.0₃0864'
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.36'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.64%  '
$ws.Range("E30").Value = '  +1.05%  '
$ws.Range("E31").Value = '  -0.59%  '
$ws.Range("E32").Value = '  +1.06%  '
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.32'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.85'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '25.06'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.66%  '
$ws.Range("D38").Value = '2.687.59'
$ws.Range("E38").Value = '  +7.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0694'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.69'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.29'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.723'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.16'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0290'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.36%  '
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = '3.206.42'
$ws.Range("E46").Value = '  +1.48%  '
$ws.Range("E47").Value = '  +13.97%  '
$ws.Range("E48").Value = '  +2.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.979'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.749'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.78%  '
